$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")
$ws.Activate()

# --- Progress (Avancement) values updated ---
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 0.5
$ws.Range("C11").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("C14").Value = 1
$ws.Range("C17").Value = 0.5
$ws.Range("C20").Value = 0.3
$ws.Range("C21").Value = 0.3

# --- Comment (D column) updates ---
# D10 had "?" - clear it out
$ws.Range("D10").ClearContents()

# D18's comment text is updated first so the shared-string table gets the
# new strings in the same order as the authored workbook (D18 then D17).
$ws.Range("D18").Value = "Méthode du coude silhouette sample"

# D17 previously had no comment; it now gets a new note about next steps
$ws.Range("D17").Value = "Mettre  plus de méthode sur les étape de clustering et le faire avec le DBScan"

# --- Selection / view state ---
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D17").Select()
